$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "GAZ"
$ws.Range("C19").Value = "entity [BFO:0000001]"
$ws.Range("D19").Value = "geographic location [GAZ:00000448]"
$ws.Range("E19").Value = "all"
